# Scraper data + show data MVC
# A new tender ("Mise en Conformité des lignes électriques HTA et BT...")
# was scraped and needs to be inserted as the 3rd data row (row 3), pushing
# the existing rows down by one. Its date_limite is unknown, so "N/A" is
# used like for every other row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new tender row, shifting rows 3..6 down to 4..7.
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "Mise en Conformité des lignes électriques HTA et BT interceptés par le projet de l’autoroute Rabat-Casablanca Continentale"
$ws.Range("B3").Value = "N/A"
